$wb = $excel.ActiveWorkbook

# --- ItemData sheet: add Alias[...] annotations to the PairField/PairFieldList
# type-definition row (row 2), matching the "Add DefaultAttr & AliasAttr" commit.
$wsItem = $wb.Worksheets.Item("ItemData")

$wsItem.Range("G2").Value  = "int_bool" + [char]10 + "Alias[PairItemIntBool]"
$wsItem.Range("I2").Value  = "int_int" + [char]10 + "Alias[PairItemIntInt64]"
$wsItem.Range("J2").Value  = "list_int_int64" + [char]10 + "Alias[PairItemIntInt64]"
$wsItem.Range("K2").Value  = "list_int_int64" + [char]10 + "Alias[PairItemIntInt64]"

# --- Update stored cursor/selection on each sheet to match the saved workbook
$wsEnumType = $wb.Worksheets.Item("@EnumConfig_ItemType")
$wsEnumType.Activate() | Out-Null
$wsEnumType.Range("C5").Select() | Out-Null

$wsEnumSubType = $wb.Worksheets.Item("@EnumConfig_ItemSubType")
$wsEnumSubType.Activate() | Out-Null
$wsEnumSubType.Range("C33").Select() | Out-Null

# ItemData keeps the selection/active-tab state, so activate & select it last
$wsItem.Activate() | Out-Null
$wsItem.Range("I19").Select() | Out-Null
